# Applies updated 'want-to-go' (F) and 'min ticket price' (G) counts
# scraped figures, per commit 456a3b4 (gh-pages data refresh).
$wb = $excel.ActiveWorkbook
$mismatches = New-Object System.Collections.ArrayList

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$cell = $ws.Range("F2")
if ($cell.Value2 -ne 552) { [void]$mismatches.Add("Sheet 1 F2: expected 552, found " + $cell.Value2) }
$cell.Value = 553
$cell = $ws.Range("F4")
if ($cell.Value2 -ne 6020) { [void]$mismatches.Add("Sheet 1 F4: expected 6020, found " + $cell.Value2) }
$cell.Value = 6024
$cell = $ws.Range("G4")
if ($cell.Value2 -ne 80) { [void]$mismatches.Add("Sheet 1 G4: expected 80, found " + $cell.Value2) }
$cell.Value = 90
$cell = $ws.Range("F12")
if ($cell.Value2 -ne 694) { [void]$mismatches.Add("Sheet 1 F12: expected 694, found " + $cell.Value2) }
$cell.Value = 696
$cell = $ws.Range("F14")
if ($cell.Value2 -ne 4) { [void]$mismatches.Add("Sheet 1 F14: expected 4, found " + $cell.Value2) }
$cell.Value = 6
$cell = $ws.Range("F15")
if ($cell.Value2 -ne 1660) { [void]$mismatches.Add("Sheet 1 F15: expected 1660, found " + $cell.Value2) }
$cell.Value = 1663
$cell = $ws.Range("F17")
if ($cell.Value2 -ne 201) { [void]$mismatches.Add("Sheet 1 F17: expected 201, found " + $cell.Value2) }
$cell.Value = 204
$cell = $ws.Range("F18")
if ($cell.Value2 -ne 657) { [void]$mismatches.Add("Sheet 1 F18: expected 657, found " + $cell.Value2) }
$cell.Value = 659
$cell = $ws.Range("F19")
if ($cell.Value2 -ne 4695) { [void]$mismatches.Add("Sheet 1 F19: expected 4695, found " + $cell.Value2) }
$cell.Value = 4703
$cell = $ws.Range("F20")
if ($cell.Value2 -ne 115) { [void]$mismatches.Add("Sheet 1 F20: expected 115, found " + $cell.Value2) }
$cell.Value = 116
$cell = $ws.Range("F22")
if ($cell.Value2 -ne 678) { [void]$mismatches.Add("Sheet 1 F22: expected 678, found " + $cell.Value2) }
$cell.Value = 679
$cell = $ws.Range("F27")
if ($cell.Value2 -ne 19) { [void]$mismatches.Add("Sheet 1 F27: expected 19, found " + $cell.Value2) }
$cell.Value = 20
$cell = $ws.Range("F29")
if ($cell.Value2 -ne 49) { [void]$mismatches.Add("Sheet 1 F29: expected 49, found " + $cell.Value2) }
$cell.Value = 50
$cell = $ws.Range("F30")
if ($cell.Value2 -ne 344) { [void]$mismatches.Add("Sheet 1 F30: expected 344, found " + $cell.Value2) }
$cell.Value = 345
$cell = $ws.Range("F31")
if ($cell.Value2 -ne 5) { [void]$mismatches.Add("Sheet 1 F31: expected 5, found " + $cell.Value2) }
$cell.Value = 6
$cell = $ws.Range("F33")
if ($cell.Value2 -ne 465) { [void]$mismatches.Add("Sheet 1 F33: expected 465, found " + $cell.Value2) }
$cell.Value = 466
$cell = $ws.Range("F36")
if ($cell.Value2 -ne 27) { [void]$mismatches.Add("Sheet 1 F36: expected 27, found " + $cell.Value2) }
$cell.Value = 28
$cell = $ws.Range("F39")
if ($cell.Value2 -ne 1287) { [void]$mismatches.Add("Sheet 1 F39: expected 1287, found " + $cell.Value2) }
$cell.Value = 1289
$cell = $ws.Range("F40")
if ($cell.Value2 -ne 1273) { [void]$mismatches.Add("Sheet 1 F40: expected 1273, found " + $cell.Value2) }
$cell.Value = 1277

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$cell = $ws.Range("F13")
if ($cell.Value2 -ne 110) { [void]$mismatches.Add("Sheet 2 F13: expected 110, found " + $cell.Value2) }
$cell.Value = 111
$cell = $ws.Range("F14")
if ($cell.Value2 -ne 215) { [void]$mismatches.Add("Sheet 2 F14: expected 215, found " + $cell.Value2) }
$cell.Value = 216
$cell = $ws.Range("F15")
if ($cell.Value2 -ne 70) { [void]$mismatches.Add("Sheet 2 F15: expected 70, found " + $cell.Value2) }
$cell.Value = 75

# --- Sheet 3: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item(3)
$cell = $ws.Range("F3")
if ($cell.Value2 -ne 767) { [void]$mismatches.Add("Sheet 3 F3: expected 767, found " + $cell.Value2) }
$cell.Value = 771
$cell = $ws.Range("F4")
if ($cell.Value2 -ne 214) { [void]$mismatches.Add("Sheet 3 F4: expected 214, found " + $cell.Value2) }
$cell.Value = 216

# --- Sheet 4: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item(4)
$cell = $ws.Range("F3")
if ($cell.Value2 -ne 552) { [void]$mismatches.Add("Sheet 4 F3: expected 552, found " + $cell.Value2) }
$cell.Value = 553
$cell = $ws.Range("F6")
if ($cell.Value2 -ne 767) { [void]$mismatches.Add("Sheet 4 F6: expected 767, found " + $cell.Value2) }
$cell.Value = 771
$cell = $ws.Range("F7")
if ($cell.Value2 -ne 214) { [void]$mismatches.Add("Sheet 4 F7: expected 214, found " + $cell.Value2) }
$cell.Value = 216
$cell = $ws.Range("F8")
if ($cell.Value2 -ne 6020) { [void]$mismatches.Add("Sheet 4 F8: expected 6020, found " + $cell.Value2) }
$cell.Value = 6024
$cell = $ws.Range("G8")
if ($cell.Value2 -ne 80) { [void]$mismatches.Add("Sheet 4 G8: expected 80, found " + $cell.Value2) }
$cell.Value = 90
$cell = $ws.Range("F23")
if ($cell.Value2 -ne 110) { [void]$mismatches.Add("Sheet 4 F23: expected 110, found " + $cell.Value2) }
$cell.Value = 111
$cell = $ws.Range("F24")
if ($cell.Value2 -ne 1660) { [void]$mismatches.Add("Sheet 4 F24: expected 1660, found " + $cell.Value2) }
$cell.Value = 1663
$cell = $ws.Range("F26")
if ($cell.Value2 -ne 201) { [void]$mismatches.Add("Sheet 4 F26: expected 201, found " + $cell.Value2) }
$cell.Value = 204
$cell = $ws.Range("F27")
if ($cell.Value2 -ne 657) { [void]$mismatches.Add("Sheet 4 F27: expected 657, found " + $cell.Value2) }
$cell.Value = 659
$cell = $ws.Range("F28")
if ($cell.Value2 -ne 4695) { [void]$mismatches.Add("Sheet 4 F28: expected 4695, found " + $cell.Value2) }
$cell.Value = 4703
$cell = $ws.Range("F35")
if ($cell.Value2 -ne 19) { [void]$mismatches.Add("Sheet 4 F35: expected 19, found " + $cell.Value2) }
$cell.Value = 20
$cell = $ws.Range("F37")
if ($cell.Value2 -ne 49) { [void]$mismatches.Add("Sheet 4 F37: expected 49, found " + $cell.Value2) }
$cell.Value = 50
$cell = $ws.Range("F38")
if ($cell.Value2 -ne 344) { [void]$mismatches.Add("Sheet 4 F38: expected 344, found " + $cell.Value2) }
$cell.Value = 345
$cell = $ws.Range("F40")
if ($cell.Value2 -ne 465) { [void]$mismatches.Add("Sheet 4 F40: expected 465, found " + $cell.Value2) }
$cell.Value = 466
$cell = $ws.Range("F45")
if ($cell.Value2 -ne 27) { [void]$mismatches.Add("Sheet 4 F45: expected 27, found " + $cell.Value2) }
$cell.Value = 28

if ($mismatches.Count -gt 0) { Write-Output ("MISMATCHES: " + ($mismatches -join "; ")) } else { Write-Output "All source values matched expectations; updates applied." }
